$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.073.85"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.790.12"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.297"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.23%  "
$ws.Range("E10").Value = "  -3.36%  "
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "2.047.47"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("E13").Value = "  +4.96%  "
$ws.Range("D14").Value = "1.791.71"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "34.056.36"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("E24").Value = "  -2.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  +1.30%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +3.15%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("D35").Value = "1.408.09"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.656"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("E38").Value = "  +2.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.923"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +13.99%  "
$ws.Range("E45").Value = "  +5.05%  "
$ws.Range("E46").Value = "  +2.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0507"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "1.949.26"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  -0.03%  "
